$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 584.3333
$ws.Range("I5").Value = 327.75
$ws.Range("K5").Value = 327.75
$ws.Range("M5").Value = -212.75
$ws.Range("H53").Value = 495.7
$ws.Range("I53").Value = 255.8
$ws.Range("K53").Value = 255.8
$ws.Range("M53").Value = 381.2
$ws.Range("H64").Value = 3916.7
$ws.Range("I64").Value = 3458.5
$ws.Range("K64").Value = 3458.5
$ws.Range("M64").Value = -3210.5
$ws.Range("H67").Value = 3916.7
$ws.Range("I67").Value = 3458.5
$ws.Range("K67").Value = 3458.5
$ws.Range("M67").Value = -2600.5
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H80").Value = 456055.62
$ws.Range("I80").Value = 747.7143
$ws.Range("J80").Value = 668532.7
$ws.Range("K80").Value = 2243.1429
$ws.Range("L80").Value = 2005598.1
$ws.Range("M80").Value = -1245.1429
$ws.Range("N80").Value = -2007594.1
$ws.Range("H83").Value = 456055.62
$ws.Range("I83").Value = 747.7143
$ws.Range("J83").Value = 668532.7
$ws.Range("K83").Value = 6729.428699999999
$ws.Range("L83").Value = 6016794.3
$ws.Range("M83").Value = -1737.428699999999
$ws.Range("N83").Value = -6026778.3
$ws.Range("H87").Value = 48607.125
$ws.Range("J87").Value = 48607.125
$ws.Range("L87").Value = 48607.125
$ws.Range("N87").Value = -51103.125
$ws.Range("H90").Value = 48607.125
$ws.Range("J90").Value = 48607.125
$ws.Range("L90").Value = 145821.375
$ws.Range("N90").Value = -158301.375
$ws.Range("H129").Value = 4238.1763
$ws.Range("I129").Value = 975.1
$ws.Range("K129").Value = 2925.3
$ws.Range("M129").Value = 2074.7
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 312.25
$ws.Range("I4").Value = 373.5
$ws.Range("J4").Value = 251
$ws.Range("K4").Value = 373.5
$ws.Range("L4").Value = 251
$ws.Range("M4").Value = -257.5
$ws.Range("N4").Value = -483
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 52645.332
$ws.Range("J82").Value = 70052.39999999999
$ws.Range("L82").Value = 70052.39999999999
$ws.Range("N82").Value = -70818.39999999999
$ws.Range("H85").Value = 52645.332
$ws.Range("J85").Value = 70052.39999999999
$ws.Range("L85").Value = 70052.39999999999
$ws.Range("N85").Value = -72704.39999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 95.625
$ws.Range("I7").Value = 103.07143
$ws.Range("K7").Value = 103.07143
$ws.Range("M7").Value = 9.928569999999993
$ws.Range("H23").Value = 27218.5
$ws.Range("I23").Value = 9900
$ws.Range("K23").Value = 9900
$ws.Range("M23").Value = -9660
$ws.Range("H27").Value = 27218.5
$ws.Range("I27").Value = 9900
$ws.Range("K27").Value = 9900
$ws.Range("M27").Value = -9708
$ws.Range("H74").Value = 57371.375
$ws.Range("J74").Value = 57371.375
$ws.Range("L74").Value = 57371.375
$ws.Range("N74").Value = -59119.375
$ws.Range("H77").Value = 57371.375
$ws.Range("J77").Value = 57371.375
$ws.Range("L77").Value = 172114.125
$ws.Range("N77").Value = -180850.125
$ws.Range("H107").Value = 502.44446
$ws.Range("I107").Value = 364.66666
$ws.Range("K107").Value = 364.66666
$ws.Range("M107").Value = 1555.33334
$ws.Range("H134").Value = 4756.32
$ws.Range("I134").Value = 3968.7
$ws.Range("K134").Value = 11906.1
$ws.Range("M134").Value = -9371.099999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3316.1667
$ws.Range("I3").Value = 2985.4
$ws.Range("J3").Value = 4970
$ws.Range("K3").Value = 8956.200000000001
$ws.Range("L3").Value = 14910
$ws.Range("M3").Value = -8844.200000000001
$ws.Range("N3").Value = -15134
$ws.Range("H11").Value = 362.04
$ws.Range("I11").Value = 391.5
$ws.Range("J11").Value = 244.2
$ws.Range("K11").Value = 1174.5
$ws.Range("L11").Value = 732.5999999999999
$ws.Range("M11").Value = -1034.5
$ws.Range("N11").Value = -1012.6
$ws.Range("H20").Value = 1704.75
$ws.Range("J20").Value = 1973
$ws.Range("L20").Value = 5919
$ws.Range("N20").Value = -6373
$ws.Range("H51").Value = 3983.5
$ws.Range("J51").Value = 5075.5
$ws.Range("L51").Value = 15226.5
$ws.Range("N51").Value = -16146.5
$ws.Range("H57").Value = 1796.3334
$ws.Range("I57").Value = 1796.3334
$ws.Range("K57").Value = 5389.0002
$ws.Range("M57").Value = -4830.0002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 868533.5
$ws.Range("I11").Value = 502500.5
$ws.Range("J11").Value = 1234566.5
$ws.Range("K11").Value = 502500.5
$ws.Range("L11").Value = 1234566.5
$ws.Range("M11").Value = -502361.5
$ws.Range("N11").Value = -1234844.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 211.3125
$ws.Range("I55").Value = 133.61539
$ws.Range("K55").Value = 133.61539
$ws.Range("M55").Value = 39.38461000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 36527.5
$ws.Range("I49").Value = 36527.5
$ws.Range("K49").Value = 36527.5
$ws.Range("M49").Value = -36297.5
$ws.Range("H62").Value = 3971.1667
$ws.Range("I62").Value = 3999.0667
$ws.Range("J62").Value = 3831.6667
$ws.Range("K62").Value = 3999.0667
$ws.Range("L62").Value = 3831.6667
$ws.Range("M62").Value = -3375.0667
$ws.Range("N62").Value = -5079.6667
$ws.Range("H65").Value = 3971.1667
$ws.Range("I65").Value = 3999.0667
$ws.Range("J65").Value = 3831.6667
$ws.Range("K65").Value = 19995.3335
$ws.Range("L65").Value = 19158.3335
$ws.Range("M65").Value = -16875.3335
$ws.Range("N65").Value = -25398.3335
$ws.Range("H81").Value = 2998.25
$ws.Range("I81").Value = 2329.6667
$ws.Range("J81").Value = 3399.4
$ws.Range("K81").Value = 4659.3334
$ws.Range("L81").Value = 6798.8
$ws.Range("M81").Value = -3598.3334
$ws.Range("N81").Value = -8920.799999999999
$ws.Range("H84").Value = 2998.25
$ws.Range("I84").Value = 2329.6667
$ws.Range("J84").Value = 3399.4
$ws.Range("K84").Value = 23296.667
$ws.Range("L84").Value = 33994
$ws.Range("M84").Value = -17992.667
$ws.Range("N84").Value = -44602
$ws.Range("H136").Value = 5130947
$ws.Range("J136").Value = 2249.8845
$ws.Range("L136").Value = 6749.6535
$ws.Range("N136").Value = -11849.6535
